$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4457773403902934
$ws.Range("C2").Value = 0.03961205890428232
$ws.Range("D2").Value = 0.654004891235644
$ws.Range("E2").Value = 0.2669223542203767
$ws.Range("G2").Value = 0.8049703976024176
$ws.Range("H2").Value = 0.8844547748153957
$ws.Range("J2").Value = 0.1390707372402176
$ws.Range("K2").Value = 0.4086094500595436
$ws.Range("N2").Value = 1.609173688003291
$ws.Range("O2").Value = 3.397180040785344
$ws.Range("B3").Value = 0.4076390570984643
$ws.Range("C3").Value = 0.03459744759958028
$ws.Range("D3").Value = 0.6434280814394526
$ws.Range("E3").Value = 0.2617436079925355
$ws.Range("G3").Value = 0.8067435465804209
$ws.Range("H3").Value = 0.8892914023486611
$ws.Range("J3").Value = 0.1356459366700093
$ws.Range("K3").Value = 0.3681200683085706
$ws.Range("N3").Value = 1.626051599796278
$ws.Range("O3").Value = 3.410777893515956
$ws.Range("B4").Value = 0.384328612579651
$ws.Range("C4").Value = 0.03150579825226885
$ws.Range("D4").Value = 0.6372624722856983
$ws.Range("E4").Value = 0.2587040662743121
$ws.Range("G4").Value = 0.8083018231257384
$ws.Range("H4").Value = 0.8926161774793258
$ws.Range("J4").Value = 0.1336204105440046
$ws.Range("K4").Value = 0.3433127508964731
$ws.Range("N4").Value = 1.636944645415465
$ws.Range("O4").Value = 3.420855622252063
$ws.Range("B5").Value = 0.3748567659484934
$ws.Range("C5").Value = 0.03024280873681562
$ws.Range("D5").Value = 0.6348327193085765
$ws.Range("E5").Value = 0.2575007389106929
$ws.Range("G5").Value = 0.8090548698592741
$ws.Range("H5").Value = 0.8940604069534217
$ws.Range("J5").Value = 0.1328144454321318
$ws.Range("K5").Value = 0.3332175297676656
$ws.Range("N5").Value = 1.641516908130357
$ws.Range("O5").Value = 3.425397112419617
$ws.Range("B6").Value = 0.3732856407191605
$ws.Range("C6").Value = 0.03003290394406122
$ws.Range("D6").Value = 0.6344342651258614
$ws.Range("E6").Value = 0.2573030614297949
$ws.Range("G6").Value = 0.8091870407254049
$ws.Range("H6").Value = 0.8943056192136751
$ws.Range("J6").Value = 0.1326817908148357
$ws.Range("K6").Value = 0.3315420855525701
$ws.Range("N6").Value = 1.642284177399021
$ws.Range("O6").Value = 3.426177482239282
$ws.Range("B7").Value = 0.3842007604943944
$ws.Range("C7").Value = 0.03148877766807345
$ws.Range("D7").Value = 0.6372293683959356
$ws.Range("E7").Value = 0.2586876947443884
$ws.Range("G7").Value = 0.8083115011030202
$ws.Range("H7").Value = 0.8926352929662329
$ws.Range("J7").Value = 0.133609462235377
$ws.Range("K7").Value = 0.3431765459491771
$ws.Range("N7").Value = 1.637005768991553
$ws.Range("O7").Value = 3.420915110072812
$ws.Range("B8").Value = 0.4326054436344862
$ws.Range("C8").Value = 0.03788568588953467
$ws.Range("D8").Value = 0.6502898955052387
$ws.Range("E8").Value = 0.2651076449146572
$ws.Range("G8").Value = 0.8054842987844779
$ws.Range("H8").Value = 0.8860487927834413
$ws.Range("J8").Value = 0.137873824839815
$ws.Range("K8").Value = 0.3946379621615677
$ws.Range("N8").Value = 1.614883204838486
$ws.Range("O8").Value = 3.401509845013692
$ws.Range("B9").Value = 0.5283539917202518
$ws.Range("C9").Value = 0.0503272690927048
$ws.Range("D9").Value = 0.6785040599080503
$ws.Range("E9").Value = 0.2788086304863171
$ws.Range("G9").Value = 0.8036684679339743
$ws.Range("H9").Value = 0.8759470811972818
$ws.Range("J9").Value = 0.146849782313673
$ws.Range("K9").Value = 0.4959575139530727
$ws.Range("N9").Value = 1.575705082843606
$ws.Range("O9").Value = 3.377171528801171
$ws.Range("B10").Value = 0.5991865361275757
$ws.Range("C10").Value = 0.05940323406045422
$ws.Range("D10").Value = 0.7008162487536254
$ws.Range("E10").Value = 0.2895521863037871
$ws.Range("G10").Value = 0.8046121460382238
$ws.Range("H10").Value = 0.8702376075742961
$ws.Range("J10").Value = 0.1538194818638772
$ws.Range("K10").Value = 0.5706247427923756
$ws.Range("N10").Value = 1.549482023029506
$ws.Range("O10").Value = 3.367654958160983
$ws.Range("B11").Value = 0.6315122225988716
$ws.Range("C11").Value = 0.06351761795792754
$ws.Range("D11").Value = 0.7113099119124797
$ws.Range("E11").Value = 0.2945868642447778
$ws.Range("G11").Value = 0.8055371886130729
$ws.Range("H11").Value = 0.8680113285198701
$ws.Range("J11").Value = 0.1570718919641791
$ws.Range("K11").Value = 0.6046387818445567
$ws.Range("N11").Value = 1.538108208366505
$ws.Range("O11").Value = 3.365143029467276
$ws.Range("B12").Value = 0.6437675645760521
$ws.Range("C12").Value = 0.06507351382190052
$ws.Range("D12").Value = 0.7153329073998407
$ws.Range("E12").Value = 0.2965145314945801
$ws.Range("G12").Value = 0.8059588376664379
$ws.Range("H12").Value = 0.8672215830888206
$ws.Range("J12").Value = 0.1583152684353877
$ws.Range("K12").Value = 0.6175253866161938
$ws.Range("N12").Value = 1.533881084910196
$ws.Range("O12").Value = 3.364453171136802
$ws.Range("B13").Value = 0.6411275305748063
$ws.Range("C13").Value = 0.06473851979450274
$ws.Range("D13").Value = 0.7144642944005
$ws.Range("E13").Value = 0.2960984342643798
$ws.Range("G13").Value = 0.8058648533952066
$ws.Range("H13").Value = 0.8673892992340484
$ws.Range("J13").Value = 0.1580469622567904
$ws.Range("K13").Value = 0.6147497585025121
$ws.Range("N13").Value = 1.534787917346793
$ws.Range("O13").Value = 3.364590119977635
$ws.Range("B14").Value = 0.6325201936388396
$ws.Range("C14").Value = 0.06364566559433626
$ws.Range("D14").Value = 0.7116398995633801
$ws.Range("E14").Value = 0.294745031115518
$ws.Range("G14").Value = 0.805570447487014
$ws.Range("H14").Value = 0.8679452877965588
$ws.Range("J14").Value = 0.1571739496886835
$ws.Range("K14").Value = 0.6056988492837263
$ws.Range("N14").Value = 1.537758838851849
$ws.Range("O14").Value = 3.365081036228844
$ws.Range("B15").Value = 0.627249796717507
$ws.Range("C15").Value = 0.06297598102074176
$ws.Range("D15").Value = 0.7099162883334316
$ws.Range("E15").Value = 0.2939187847135258
$ws.Range("G15").Value = 0.8053994099142443
$ws.Range("H15").Value = 0.8682927863702048
$ws.Range("J15").Value = 0.1566407353419521
$ws.Range("K15").Value = 0.60015569879792
$ws.Range("N15").Value = 1.539589020298434
$ws.Range("O15").Value = 3.365415773948911
$ws.Range("B16").Value = 0.5970760106875161
$ws.Range("C16").Value = 0.05913405502705871
$ws.Range("D16").Value = 0.7001373667320081
$ws.Range("E16").Value = 0.2892261197385508
$ws.Range("G16").Value = 0.8045616739127865
$ws.Range("H16").Value = 0.8703905611388194
$ws.Range("J16").Value = 0.153608575416925
$ws.Range("K16").Value = 0.5684027498625142
$ws.Range("N16").Value = 1.550236505247844
$ws.Range("O16").Value = 3.367855689074474
$ws.Range("B17").Value = 0.578591485000203
$ws.Range("C17").Value = 0.05677343862241457
$ws.Range("D17").Value = 0.6942262435043745
$ws.Range("E17").Value = 0.2863850345919943
$ws.Range("G17").Value = 0.8041747727736492
$ws.Range("H17").Value = 0.871772461727943
$ws.Range("J17").Value = 0.1517694011775035
$ws.Range("K17").Value = 0.5489350867579788
$ws.Range("N17").Value = 1.556910629501633
$ws.Range("O17").Value = 3.36981795833421
$ws.Range("B18").Value = 0.5679694621100566
$ws.Range("C18").Value = 0.05541433118008854
$ws.Range("D18").Value = 0.6908586852394194
$ws.Range("E18").Value = 0.284764792674288
$ws.Range("G18").Value = 0.8039988966946083
$ws.Range("H18").Value = 0.8726022165791818
$ws.Range("J18").Value = 0.1507192624658273
$ws.Range("K18").Value = 0.5377423132054275
$ws.Range("N18").Value = 1.560801662390318
$ws.Range("O18").Value = 3.371117645259972
$ws.Range("B19").Value = 0.5643747267124013
$ws.Range("C19").Value = 0.05495393226495082
$ws.Range("D19").Value = 0.6897240510715505
$ws.Range("E19").Value = 0.2842185911564243
$ws.Range("G19").Value = 0.8039473598277453
$ws.Range("H19").Value = 0.8728891567571537
$ws.Range("J19").Value = 0.1503650273530326
$ws.Range("K19").Value = 0.5339534258203855
$ws.Range("N19").Value = 1.562128072075476
$ws.Range("O19").Value = 3.371587072195069
$ws.Range("B20").Value = 0.5805581856481581
$ws.Range("C20").Value = 0.05702486982777089
$ws.Range("D20").Value = 0.694852144293975
$ws.Range("E20").Value = 0.2866860372891864
$ws.Range("G20").Value = 0.8042111295256973
$ws.Range("H20").Value = 0.8716217420964512
$ws.Range("J20").Value = 0.1519643869931286
$ws.Range("K20").Value = 0.5510069902394434
$ws.Range("N20").Value = 1.556194749162461
$ws.Range("O20").Value = 3.369591368952342
$ws.Range("B21").Value = 0.6350479936018303
$ws.Range("C21").Value = 0.06396672180905227
$ws.Range("D21").Value = 0.7124681567821654
$ws.Range("E21").Value = 0.2951419849793666
$ws.Range("G21").Value = 0.8056549845465355
$ws.Range("H21").Value = 0.8677805343755978
$ws.Range("J21").Value = 0.1574300555922292
$ws.Range("K21").Value = 0.6083571565991122
$ws.Range("N21").Value = 1.536884038430549
$ws.Range("O21").Value = 3.364929748798744
$ws.Range("B22").Value = 0.6707433446387085
$ws.Range("C22").Value = 0.06849116244715958
$ws.Range("D22").Value = 0.7242683672116073
$ws.Range("E22").Value = 0.3007916558517252
$ws.Range("G22").Value = 0.8070145721653148
$ws.Range("H22").Value = 0.8655807227105043
$ws.Range("J22").Value = 0.1610707228019237
$ws.Range("K22").Value = 0.6458748145659854
$ws.Range("N22").Value = 1.5247290810044
$ws.Range("O22").Value = 3.36340650482029
$ws.Range("B23").Value = 0.6516846551118931
$ws.Range("C23").Value = 0.06607754763186335
$ws.Range("D23").Value = 0.7179441522688421
$ws.Range("E23").Value = 0.2977650618535108
$ws.Range("G23").Value = 0.8062508539558166
$ws.Range("H23").Value = 0.8667263963176168
$ws.Range("J23").Value = 0.1591213630534583
$ws.Range("K23").Value = 0.6258478412107991
$ws.Range("N23").Value = 1.531173778254799
$ws.Range("O23").Value = 3.364080084023101
$ws.Range("B24").Value = 0.5796690247369725
$ws.Range("C24").Value = 0.05691120387331239
$ws.Range("D24").Value = 0.694569078519379
$ws.Range("E24").Value = 0.2865499130507203
$ws.Range("G24").Value = 0.8041945476152534
$ws.Range("H24").Value = 0.8716897725354045
$ws.Range("J24").Value = 0.1518762113943808
$ws.Range("K24").Value = 0.5500702843576448
$ws.Range("N24").Value = 1.556518230489425
$ws.Range("O24").Value = 3.369693275635058
$ws.Range("B25").Value = 0.5023647998183947
$ws.Range("C25").Value = 0.04697272120016294
$ws.Range("D25").Value = 0.6705931054808332
$ws.Range("E25").Value = 0.2749831732047596
$ws.Range("G25").Value = 0.8037600512547129
$ws.Range("H25").Value = 0.8783789232575145
$ws.Range("J25").Value = 0.144355755982879
$ws.Range("K25").Value = 0.4685065618421902
$ws.Range("N25").Value = 1.585853922645285
$ws.Range("O25").Value = 3.382287003045036
